# Generate Report for Handback
# Row 7 ("23b6880d-...") previously shared its Correspond Handback DateTime
# values with row 8 ("2b8b8f24-..."). This gives row 7 its own, distinct
# handback timestamps on both the zh-cn and de-de sheets while row 8 keeps
# its original values.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D7").Value = "2016-02-29 04:01:06"
$wsZhCn.Range("G7").Value = "2016-02-29 04:02:06"
$wsZhCn.Range("D8").Value = "2016-02-29 03:58:46"
$wsZhCn.Range("G8").Value = "2016-02-29 03:59:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D7").Value = "2016-02-29 04:01:21"
$wsDeDe.Range("G7").Value = "2016-02-29 04:02:29"
$wsDeDe.Range("D8").Value = "2016-02-29 03:58:57"
$wsDeDe.Range("G8").Value = "2016-02-29 04:00:15"
